# Build site at 2021-10-27 12:19:33 UTC
# Applies the LOQ4062.xlsx changes:
#  - "Semestre ideal" value updated from EB-8,EQD-8,EQN-10 to EB-8,EQD-7,EQN-9
#  - "Requisitos" list loses the "LOQ4057 - Operacoes Unitarias III" weak-requisite
#    entry; the remaining two requisite rows shift up by one, and the now-empty
#    trailing row is removed from the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Semestre ideal:" value in B9/C9
$ws.Range("B9").Value = "EB-8,EQD-7,EQN-9"
$ws.Range("C9").Value = "EB-8,EQD-7,EQN-9"

# 2. Shift the two remaining "Requisitos" rows up, dropping the LOQ4057 entry.
$req2 = "LOQ4085 -  Operações Unitárias I  (Requisito fraco)`n"
$req3 = "LOQ4086 -  Operações Unitárias II  (Requisito fraco)`n"

$ws.Range("B25").Value = $req2
$ws.Range("C25").Value = $req2

$ws.Range("B26").Value = $req3
$ws.Range("C26").Value = $req3

# 3. Remove the now-redundant last row (former row 27) entirely, shifting the
#    dimension from A1:C27 down to A1:C26.
$ws.Rows(27).Delete()
